$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D3").Value = -7.039899999999994
$ws.Range("C7").Value = -13.19759999999999
$ws.Range("A10").Value = -21.83119999999998
$ws.Range("E10").Value = 16.4128
$ws.Range("A12").Value = -21.5407
$ws.Range("E14").Value = 16.68210000000001
$ws.Range("C15").Value = -14.35129999999998
$ws.Range("A18").Value = -22.06490000000001
$ws.Range("D18").Value = -8.099199999999993
$ws.Range("D19").Value = -9.043099999999992
$ws.Range("C20").Value = -12.0858
$ws.Range("D27").Value = -9.047600000000006
$ws.Range("C29").Value = -11.58900000000001
$ws.Range("C30").Value = -12.75639999999999
$ws.Range("C31").Value = -13.0253
$ws.Range("E32").Value = 16.65799999999999
$ws.Range("E35").Value = 16.641
$ws.Range("A37").Value = -20.68060000000001
$ws.Range("C40").Value = -13.03000000000001
$ws.Range("D42").Value = -8.779599999999993
$ws.Range("E43").Value = 17.3561
$ws.Range("D44").Value = -7.761500000000001
$ws.Range("D47").Value = -7.703000000000003
$ws.Range("E49").Value = 15.2932
$ws.Range("A55").Value = -22.33039999999999
$ws.Range("E56").Value = 16.5475
$ws.Range("D58").Value = -8.328699999999991
$ws.Range("A68").Value = -21.6879
$ws.Range("C68").Value = -12.0696
$ws.Range("E69").Value = 17.46070000000002
$ws.Range("D73").Value = -7.952899999999997
$ws.Range("C76").Value = -11.89210000000001
$ws.Range("A77").Value = -20.87609999999998
$ws.Range("A78").Value = -20.65699999999998
$ws.Range("E81").Value = 16.63549999999999
$ws.Range("C87").Value = -13.69399999999998
$ws.Range("C88").Value = -13.35029999999999
$ws.Range("E92").Value = 18.43640000000001
$ws.Range("D95").Value = -7.878499999999997
$ws.Range("C96").Value = -12.93000000000001
$ws.Range("C98").Value = -12.329
$ws.Range("C101").Value = -13.5694
$ws.Range("D101").Value = -8.155299999999999
$ws.Range("C102").Value = -13.3066
